$d = $word.ActiveDocument

# Remove every paragraph except the first one (iterate backwards so
# indices of earlier paragraphs stay stable while we delete).
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $d.Paragraphs.Item($i).Range.Delete()
}

# Empty out whatever remains in the first (only) paragraph.
$d.Paragraphs.Item(1).Range.Delete()

# Insert the new paragraph content as raw OOXML so we can reproduce the
# exact run/proofErr/bookmark structure from the target document.
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t xml:space="preserve">The data is </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:t>actually normally</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:t xml:space="preserve"> distributed, but it might need transformation to reveal its normality. For example, lognormal distribution </w:t>
  </w:r>
  <w:r>
    <w:t>bec</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:r>
    <w:t xml:space="preserve">omes normal distribution after taking a log on it. </w:t>
  </w:r>
</w:p>
'@

[void]$d.Paragraphs.Item(1).Range.InsertXML($xml)
